# Weekly update: a new price record for Espinaca (Feria Lagunitas de Puerto
# Montt) is inserted as the new row 4, pushing the existing data rows
# (formerly 4-14) down to rows 5-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (shifts rows 4..14 -> 5..15,
# copying row formatting - including the date number format on column D -
# from the row above, same as Excel's native Insert behaviour).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 44659
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("N4").Value = "$/cuna 10 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 1000
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"
